$wb = $excel.ActiveWorkbook

# --- Sheet 1: Metadata ---
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/ValueSet/sex-assigned-at-birth"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Sheet 2: Include from Sex Assigned At ... ---
$inc = $wb.Worksheets.Item(2)

# System URI: ibm.com -> linuxforhealth.org
$inc.Range("B4").Value = "http://linuxforhealth.org/fhir/cdm/CodeSystem/sex-assigned-at-birth"
